$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.585.37'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.54%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.991.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.07%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4677'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.29%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3956'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.21%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.63'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08151'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.89%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.002'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.89%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.94'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.65%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.001.87'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.75%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.256'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.876'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.91%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07128'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.26%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.92'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.98%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.004'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001008'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.20%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.42'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.15%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.13%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.600.02'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.560'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.82%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.27'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.88%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.119'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.67%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.92'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.92%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.71'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.018'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.65%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '120.37'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.943'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.93%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09462'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.12%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9153'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.70%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.284'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.48%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.354'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.88%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.182'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.18%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05859'

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.177'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.61%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02128'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.78%  '

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.000003310'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +86.70%  '

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.927'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.96%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5785'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.41%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1831'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.53%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.921'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.29%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.09'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.71%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.762'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +8.41%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5405'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.34%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.217'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.79%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.874'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.91%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06970'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.04%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '114.11'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.39%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3080'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.94%  '
